$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.308.63"
$ws.Range("E2").Value = "  +2.70%  "

$ws.Range("D3").Value = "1.718.14"
$ws.Range("E3").Value = "  +2.96%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'238.74"
$ws.Range("E5").Value = "  +0.64%  "

$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").Value = "'0.4703"
$ws.Range("E7").Value = "  -2.27%  "

$ws.Range("D8").Value = "'0.2623"
$ws.Range("E8").Value = "  -0.40%  "

$ws.Range("D9").Value = "'0.06184"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("D10").Value = "1.717.72"
$ws.Range("E10").Value = "  +3.19%  "

$ws.Range("D11").Value = "'0.07068"
$ws.Range("E11").Value = "  -0.47%  "

$ws.Range("D12").Value = "'15.30"
$ws.Range("E12").Value = "  +2.99%  "

$ws.Range("D13").Value = "'0.5906"
$ws.Range("E13").Value = "  -1.67%  "

$ws.Range("D14").Value = "'4.386"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").Value = "'76.07"
$ws.Range("E15").Value = "  +1.72%  "

$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  +0.20%  "

$ws.Range("D18").Value = "26.310.36"
$ws.Range("E18").Value = "  +2.71%  "

$ws.Range("D19").Value = "'0.000006793"
$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("D20").Value = "'11.53"
$ws.Range("E20").Value = "  +0.33%  "

$ws.Range("D21").Value = "1.939.71"
$ws.Range("E21").Value = "  +3.22%  "

$ws.Range("D22").Value = "'4.538"
$ws.Range("E22").Value = "  +1.34%  "

$ws.Range("D23").Value = "'8.729"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "'5.317"
$ws.Range("E24").Value = "  -0.99%  "

$ws.Range("D25").Value = "'135.76"
$ws.Range("E25").Value = "  +0.90%  "

$ws.Range("D26").Value = "'15.23"
$ws.Range("E26").Value = "  +0.79%  "

$ws.Range("D27").Value = "'108.17"
$ws.Range("E27").Value = "  +3.04%  "

$ws.Range("D28").Value = "'1.402"
$ws.Range("E28").Value = "  -0.45%  "

$ws.Range("E29").Value = "  +3.31%  "

$ws.Range("D30").Value = "'4.000"
$ws.Range("E30").Value = "  +0.25%  "

$ws.Range("D31").Value = "'3.678"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("D32").Value = "'0.07723"
$ws.Range("E32").Value = "  +0.19%  "

$ws.Range("D33").Value = "'0.04444"
$ws.Range("E33").Value = "  +1.77%  "

$ws.Range("D34").Value = "'2.615"
$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("D35").Value = "'0.9745"
$ws.Range("E35").Value = "  +2.08%  "

$ws.Range("D36").Value = "'0.6185"
$ws.Range("E36").Value = "  +0.35%  "

$ws.Range("D37").Value = "'0.9225"
$ws.Range("E37").Value = "  +5.61%  "

$ws.Range("D38").Value = "'114.14"
$ws.Range("E38").Value = "  +16.80%  "

$ws.Range("D39").Value = "'2.414"
$ws.Range("E39").Value = "  -7.88%  "

$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("D41").Value = "'1.895"
$ws.Range("E41").Value = "  +1.15%  "

$ws.Range("D42").Value = "'0.01480"
$ws.Range("E42").Value = "  -2.56%  "

$ws.Range("D43").Value = "'5.344"
$ws.Range("E43").Value = "  +14.02%  "

$ws.Range("D44").Value = "'0.3805"
$ws.Range("E44").Value = "  +0.30%  "

$ws.Range("D45").Value = "'0.1162"
$ws.Range("E45").Value = "  +3.19%  "

$ws.Range("D46").Value = "'6.260"
$ws.Range("E46").Value = "  +0.27%  "

$ws.Range("D47").Value = "'0.05290"
$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("D48").Value = "'30.49"
$ws.Range("E48").Value = "  +3.01%  "

$ws.Range("E49").Value = "  +3.73%  "

$ws.Range("D50").Value = "'0.3372"
$ws.Range("E50").Value = "  +0.42%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.214"
$ws.Range("E51").Value = "  +1.25%  "

Write-Output "done"